# Apply odds/score updates for 2025-05-22 FlashScore weekly games sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("G7").Value = 1.9
$ws.Range("I7").Value = 3.4
$ws.Range("T7").Value = 10
$ws.Range("AD7").Value = 126

# Row 10
$ws.Range("G10").Value = 1.8
$ws.Range("I10").Value = 3.9
$ws.Range("Q10").Value = 3.25
$ws.Range("R10").Value = 1.65
$ws.Range("S10").Value = 2.12
$ws.Range("T10").Value = 9
$ws.Range("U10").Value = 9.5
$ws.Range("AB10").Value = 13
$ws.Range("AE10").Value = 13
$ws.Range("AF10").Value = 21
$ws.Range("AI10").Value = 29

# Row 11
$ws.Range("G11").Value = 1.67
$ws.Range("T11").Value = 7
$ws.Range("U11").Value = 8
$ws.Range("X11").Value = 13
$ws.Range("Y11").Value = 26
$ws.Range("Z11").Value = 10
$ws.Range("AD11").Value = 301
$ws.Range("AE11").Value = 13

# Row 13
$ws.Range("N13").Value = 1.95
$ws.Range("O13").Value = 1.85

# Row 16
$ws.Range("H16").Value = 4
$ws.Range("J16").Value = 1.02
$ws.Range("K16").Value = 12
$ws.Range("L16").Value = 1.25
$ws.Range("M16").Value = 3.75
$ws.Range("N16").Value = 1.8
$ws.Range("O16").Value = 2
$ws.Range("R16").Value = 1.83
$ws.Range("S16").Value = 1.83
$ws.Range("Z16").Value = 12
$ws.Range("AD16").Value = 600
$ws.Range("AI16").Value = 41

# Row 17
$ws.Range("G17").Value = 3.6
$ws.Range("H17").Value = 3.7
$ws.Range("I17").Value = 1.9
$ws.Range("N17").Value = 1.65
$ws.Range("O17").Value = 2.2
$ws.Range("P17").Value = 1.3
$ws.Range("Q17").Value = 3.4
$ws.Range("U17").Value = 21
$ws.Range("X17").Value = 26
$ws.Range("Y17").Value = 29
$ws.Range("AA17").Value = 7.5
$ws.Range("AD17").Value = 126
$ws.Range("AF17").Value = 11
$ws.Range("AH17").Value = 17
$ws.Range("AI17").Value = 15

# Row 18
$ws.Range("G18").Value = 1.42
$ws.Range("H18").Value = 5
$ws.Range("I18").Value = 5.5
$ws.Range("J18").Value = 21
$ws.Range("K18").Value = 1.03
$ws.Range("N18").Value = 1.44
$ws.Range("O18").Value = 2.63
$ws.Range("U18").Value = 8.5
$ws.Range("V18").Value = 8.5
$ws.Range("Z18").Value = 21
$ws.Range("AA18").Value = 10
$ws.Range("AE18").Value = 21
$ws.Range("AF18").Value = 34
$ws.Range("AG18").Value = 19
$ws.Range("AJ18").Value = 41

# Row 19
$ws.Range("J19").Value = 1.05
$ws.Range("K19").Value = 11

# Row 20
$ws.Range("G20").Value = 1.83
$ws.Range("I20").Value = 4
$ws.Range("T20").Value = 7.5
$ws.Range("AE20").Value = 11

# Row 21
$ws.Range("H21").Value = 3.7
$ws.Range("J21").Value = 1.06
$ws.Range("K21").Value = 9.5
$ws.Range("L21").Value = 1.3
$ws.Range("M21").Value = 3.4
$ws.Range("N21").Value = 2.03
$ws.Range("O21").Value = 1.78
$ws.Range("P21").Value = 1.4
$ws.Range("Q21").Value = 2.75
$ws.Range("R21").Value = 1.95
$ws.Range("S21").Value = 1.8
$ws.Range("T21").Value = 6.5
$ws.Range("Z21").Value = 9.5
$ws.Range("AB21").Value = 17
$ws.Range("AC21").Value = 51
$ws.Range("AD21").Value = 351
$ws.Range("AE21").Value = 12

# Row 26
$ws.Range("K26").Value = 17
$ws.Range("N26").Value = 1.57
$ws.Range("O26").Value = 2.35
$ws.Range("R26").Value = 1.53
$ws.Range("S26").Value = 2.38
$ws.Range("T26").Value = 10
$ws.Range("U26").Value = 11
$ws.Range("Z26").Value = 17
$ws.Range("AE26").Value = 15
$ws.Range("AF26").Value = 23
$ws.Range("AG26").Value = 13
$ws.Range("AI26").Value = 29
$ws.Range("AJ26").Value = 29

# Row 28
$ws.Range("G28").Value = 2.1
$ws.Range("H28").Value = 3.6
$ws.Range("I28").Value = 3.1
$ws.Range("L28").Value = 1.18
$ws.Range("M28").Value = 4.5
$ws.Range("T28").Value = 10
$ws.Range("U28").Value = 12
$ws.Range("V28").Value = 9
$ws.Range("W28").Value = 21
$ws.Range("X28").Value = 15
$ws.Range("AA28").Value = 7
$ws.Range("AB28").Value = 12
$ws.Range("AD28").Value = 126
$ws.Range("AF28").Value = 19
$ws.Range("AG28").Value = 11
$ws.Range("AH28").Value = 34
$ws.Range("AI28").Value = 23
$ws.Range("AJ28").Value = 26
